$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (every cell that showed the old status on all three sheets)
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2. Latest Handback DateTime (column K) for zh-cn and de-de
# ---------------------------------------------------------------------------
$zhcn.Range("K2").Value = "2016-09-06 07:10:45"
$zhcn.Range("K3").Value = "2016-09-06 07:10:45"
$dede.Range("K2").Value = "2016-09-06 07:11:10"
$dede.Range("K3").Value = "2016-09-06 07:11:10"

# ---------------------------------------------------------------------------
# 3. Latest Target File (column I, hyperlinked) and Latest Handback File
#    (column J) for zh-cn and de-de
# ---------------------------------------------------------------------------
$md97 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a7d39eb313aa6c754f8aede2ed104cf5b77d719d/e2e/97d54cc1-b806-4f0e-83ad-c1066e84df10.md"
$mda7 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a7d39eb313aa6c754f8aede2ed104cf5b77d719d/e2e/a7f655fd-4c9e-4e2e-b5bb-616325fb1b97.md"

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $md97, "", "", "97d54cc1-b806-4f0e-83ad-c1066e84df10.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $mda7, "", "", "a7f655fd-4c9e-4e2e-b5bb-616325fb1b97.md")
$zhcn.Range("J2").Value = "97d54cc1-b806-4f0e-83ad-c1066e84df10.c6567508832b715ed67937c1f38159a848ff4a53.zh-cn.xlf"
$zhcn.Range("J3").Value = "a7f655fd-4c9e-4e2e-b5bb-616325fb1b97.c22f3d74b7ee2809148bf4fdb5ad9290b81a4322.zh-cn.xlf"

$dede.Hyperlinks.Add($dede.Range("I2"), $md97, "", "", "97d54cc1-b806-4f0e-83ad-c1066e84df10.md")
$dede.Hyperlinks.Add($dede.Range("I3"), $mda7, "", "", "a7f655fd-4c9e-4e2e-b5bb-616325fb1b97.md")
$dede.Range("J2").Value = "97d54cc1-b806-4f0e-83ad-c1066e84df10.c6567508832b715ed67937c1f38159a848ff4a53.de-de.xlf"
$dede.Range("J3").Value = "a7f655fd-4c9e-4e2e-b5bb-616325fb1b97.c22f3d74b7ee2809148bf4fdb5ad9290b81a4322.de-de.xlf"

# ---------------------------------------------------------------------------
# 4. Column width changes
#    ColumnWidth is quantized by the host app to 1/6-character steps, so the
#    values below are chosen to land on (or as close as possible to) the
#    target stored widths of 29.9777050018311 / 40 characters.
# ---------------------------------------------------------------------------
$wideCol    = 30 - (5/6)       # -> stored width 30 (closest reachable to 29.9777050018311)
$fortyCol   = 40 - (5/6)       # -> stored width 40

$overview.Columns.Item(5).ColumnWidth = $wideCol
$overview.Columns.Item(6).ColumnWidth = $wideCol

$zhcn.Columns.Item(3).ColumnWidth = $wideCol
$zhcn.Columns.Item(9).ColumnWidth = $fortyCol
$zhcn.Columns.Item(10).ColumnWidth = $fortyCol

$dede.Columns.Item(3).ColumnWidth = $wideCol
$dede.Columns.Item(9).ColumnWidth = $fortyCol
$dede.Columns.Item(10).ColumnWidth = $fortyCol

Write-Host "Handback report generated"
